$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 12 data
$ws.Range("A12").Value = 41439
$ws.Range("A12").NumberFormat = 'ddd\ dd/mm/yyyy'

$ws.Range("B12").Value = 1.5
$ws.Range("C12").Value = 2.5
$ws.Range("D12").Value = "Implementation of semaphores and first, very preliminary but successfuls tests"

# Update selection to A12, matching the diff
$ws.Range("A12").Select()
